{"js": "// Replace the date line and every two-digit-by-two-digit multiplication\n// answer in the table with the new values from the commit.\nconst replacements = [\n  [\"2024-07-21 Sunday\", \"2024-07-22 Monday\"],\n  [\"71\u00d725=1775\", \"28\u00d755=1540\"],\n  [\"75\u00d738=2850\", \"43\u00d785=3655\"],\n  [\"40\u00d739=1560\", \"32\u00d741=1312\"],\n  [\"77\u00d763=4851\", \"61\u00d788=5368\"],\n  [\"89\u00d737=3293\", \"22\u00d780=1760\"],\n  [\"36\u00d722=792\", \"42\u00d741=1722\"],\n  [\"74\u00d771=5254\", \"84\u00d729=2436\"],\n  [\"52\u00d791=4732\", \"94\u00d744=4136\"],\n  [\"48\u00d722=1056\", \"18\u00d783=1494\"],\n  [\"59\u00d773=4307\", \"14\u00d723=322\"],\n  [\"44\u00d782=3608\", \"63\u00d784=5292\"],\n  [\"57\u00d785=4845\", \"51\u00d745=2295\"],\n  [\"69\u00d739=2691\", \"69\u00d735=2415\"],\n  [\"26\u00d768=1768\", \"40\u00d752=2080\"],\n  [\"12\u00d778=936\", \"21\u00d720=420\"],\n  [\"29\u00d789=2581\", \"28\u00d799=2772\"],\n  [\"32\u00d730=960\", \"15\u00d789=1335\"],\n  [\"90\u00d737=3330\", \"78\u00d797=7566\"],\n  [\"36\u00d784=3024\", \"71\u00d768=4828\"],\n  [\"63\u00d794=5922\", \"80\u00d794=7520\"],\n  [\"19\u00d718=342\", \"21\u00d790=1890\"],\n  [\"38\u00d769=2622\", \"23\u00d732=736\"],\n  [\"71\u00d798=6958\", \"43\u00d775=3225\"],\n  [\"52\u00d776=3952\", \"93\u00d749=4557\"],\n  [\"36\u00d785=3060\", \"63\u00d731=1953\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit-by-two-digit multiplication\n# answer in the table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-21 Sunday\", \"2024-07-22 Monday\"),\n    @(\"71\u00d725=1775\", \"28\u00d755=1540\"),\n    @(\"75\u00d738=2850\", \"43\u00d785=3655\"),\n    @(\"40\u00d739=1560\", \"32\u00d741=1312\"),\n    @(\"77\u00d763=4851\", \"61\u00d788=5368\"),\n    @(\"89\u00d737=3293\", \"22\u00d780=1760\"),\n    @(\"36\u00d722=792\", \"42\u00d741=1722\"),\n    @(\"74\u00d771=5254\", \"84\u00d729=2436\"),\n    @(\"52\u00d791=4732\", \"94\u00d744=4136\"),\n    @(\"48\u00d722=1056\", \"18\u00d783=1494\"),\n    @(\"59\u00d773=4307\", \"14\u00d723=322\"),\n    @(\"44\u00d782=3608\", \"63\u00d784=5292\"),\n    @(\"57\u00d785=4845\", \"51\u00d745=2295\"),\n    @(\"69\u00d739=2691\", \"69\u00d735=2415\"),\n    @(\"26\u00d768=1768\", \"40\u00d752=2080\"),\n    @(\"12\u00d778=936\", \"21\u00d720=420\"),\n    @(\"29\u00d789=2581\", \"28\u00d799=2772\"),\n    @(\"32\u00d730=960\", \"15\u00d789=1335\"),\n    @(\"90\u00d737=3330\", \"78\u00d797=7566\"),\n    @(\"36\u00d784=3024\", \"71\u00d768=4828\"),\n    @(\"63\u00d794=5922\", \"80\u00d794=7520\"),\n    @(\"19\u00d718=342\", \"21\u00d790=1890\"),\n    @(\"38\u00d769=2622\", \"23\u00d732=736\"),\n    @(\"71\u00d798=6958\", \"43\u00d775=3225\"),\n    @(\"52\u00d776=3952\", \"93\u00d749=4557\"),\n    @(\"36\u00d785=3060\", \"63\u00d731=1953\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
